$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "61.661.96"
$r.Style = "Normal"
$ws.Range("E2").Value = "  -0.60%  "
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "3.397.26"
$r.Style = "Normal"
$ws.Range("E3").Value = "  -0.72%  "
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "408.10"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -0.40%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "126.83"
$r.Style = "Normal"
$ws.Range("E6").Value = "  -1.68%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.614"
$r.Style = "Normal"
$ws.Range("E7").Value = "  -2.30%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"
$ws.Range("E8").Value = "  +0.06%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.717"
$r.Style = "Normal"
$ws.Range("E9").Value = "  -4.05%  "
$ws.Range("E10").Value = "  -8.34%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "42.16"
$r.Style = "Normal"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "3.937.56"
$r.Style = "Normal"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "9.05"
$r.Style = "Normal"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "0.140"
$r.Style = "Normal"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("E15").Value = "  -8.79%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "20.26"
$r.Style = "Normal"
$ws.Range("E16").Value = "  -3.70%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "3.401.07"
$r.Style = "Normal"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "1.07"
$r.Style = "Normal"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "12.16"
$r.Style = "Normal"
$ws.Range("E19").Value = "  -2.11%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "61.673.39"
$r.Style = "Normal"
$ws.Range("E20").Value = "  -0.60%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "483.97"
$r.Style = "Normal"
$ws.Range("E21").Value = "  +20.98%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "89.11"
$r.Style = "Normal"
$ws.Range("E22").Value = "  -0.36%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "3.20"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("E24").Value = "  -1.41%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "3.27"
$r.Style = "Normal"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("E26").Value = "  +1.44%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "9.26"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +5.84%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "4.81"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +0.18%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "7.90"
$r.Style = "Normal"
$ws.Range("E29").Value = "  +4.27%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "2.75"
$r.Style = "Normal"
$ws.Range("E30").Value = "  -0.10%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "11.74"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("E32").Value = "  -3.29%  "
$ws.Range("E33").Value = "  -6.38%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "40.79"
$r.Style = "Normal"
$ws.Range("E34").Value = "  -5.27%  "
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"
$ws.Range("E35").Value = "  -0.66%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "55.28"
$r.Style = "Normal"
$ws.Range("E36").Value = "  +2.81%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "0.0483"
$r.Style = "Normal"
$ws.Range("E37").Value = "  -3.17%  "
$ws.Range("E38").Value = "  +0.08%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.327"
$r.Style = "Normal"
$ws.Range("E39").Value = "  +4.66%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "146.96"
$r.Style = "Normal"
$ws.Range("E40").Value = "  +3.69%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "2.93"
$r.Style = "Normal"
$ws.Range("E41").Value = "  +0.53%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "3.31"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -1.73%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "0.133"
$r.Style = "Normal"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("E44").Value = "  +3.74%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "2.53"
$r.Style = "Normal"
$ws.Range("E45").Value = "  +4.60%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "4.15"
$r.Style = "Normal"
$ws.Range("E46").Value = "  +0.86%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "2.34"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +15.44%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "16.25"
$r.Style = "Normal"
$ws.Range("E48").Value = "  -2.55%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "21.84"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.143"
$r.Style = "Normal"
$ws.Range("E50").Value = "  +9.62%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "111.91"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +13.73%  "
